$d = $word.ActiveDocument

# Change company name "ОАО «НИАЭП»" -> "АО ИК «АСЭ»"
$d.Content.Find.Execute("ОАО «НИАЭП»", $true, $false, $false, $false, $false,
                         $true, 1, $false, "АО ИК «АСЭ»", 2)

# Change addressee "Начальнику ОСКР" -> "Начальнику УСКР"
$d.Content.Find.Execute("Начальнику ОСКР", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Начальнику УСКР", 2)
